# Refresh the "cryptos" price table: Price (column D) and Volume(1h)
# (column E) values for most rows, plus a Coin/Link swap + updated
# Price/Volume for rows 40-41 (Stacks <-> EthereumClassic trade places).
#
# GitHub Actions-style data refresh, mirroring the upstream commit
# "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells store prices/percentages as literal text (e.g. "555.97",
# "  +2.35%  "), not numbers. Plain `.Value = "555.97"` assignment lets
# Excel's COM layer auto-coerce digit-and-dot strings into Number cells,
# which would corrupt values like "555.97" -> 555.97000000000003 and
# flip the stored cell type. Force those through as text explicitly.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "62.244.62"
$ws.Cells.Item(2, 5).Value = "  +2.35%  "

$ws.Cells.Item(3, 4).Value = "2.419.45"
$ws.Cells.Item(3, 5).Value = "  +3.01%  "

Set-TextValue 5 4 "555.97"
$ws.Cells.Item(5, 5).Value = "  +2.19%  "

Set-TextValue 6 4 "142.74"
$ws.Cells.Item(6, 5).Value = "  +4.49%  "

$ws.Cells.Item(7, 5).Value = "  +0.00%  "

Set-TextValue 8 4 "0.533"
$ws.Cells.Item(8, 5).Value = "  +1.76%  "

$ws.Cells.Item(9, 4).Value = "2.418.33"
$ws.Cells.Item(9, 5).Value = "  +3.01%  "

Set-TextValue 10 4 "0.109"
$ws.Cells.Item(10, 5).Value = "  +4.29%  "

$ws.Cells.Item(11, 5).Value = "  -0.59%  "

$ws.Cells.Item(12, 5).Value = "  +1.48%  "

Set-TextValue 13 4 "0.351"
$ws.Cells.Item(13, 5).Value = "  +2.31%  "

Set-TextValue 14 4 "26.18"
$ws.Cells.Item(14, 5).Value = "  +6.25%  "

$ws.Cells.Item(15, 5).Value = "  +8.30%  "

$ws.Cells.Item(16, 4).Value = "2.846.35"
$ws.Cells.Item(16, 5).Value = "  +2.70%  "

$ws.Cells.Item(17, 4).Value = "62.175.13"
$ws.Cells.Item(17, 5).Value = "  +2.31%  "

$ws.Cells.Item(18, 4).Value = "2.420.82"
$ws.Cells.Item(18, 5).Value = "  +2.89%  "

Set-TextValue 19 4 "11.08"
$ws.Cells.Item(19, 5).Value = "  +4.33%  "

Set-TextValue 20 4 "4.20"

Set-TextValue 21 4 "324.02"
$ws.Cells.Item(21, 5).Value = "  +1.49%  "

Set-TextValue 22 4 "6.71"
$ws.Cells.Item(22, 5).Value = "  +2.31%  "

$ws.Cells.Item(23, 5).Value = "  +0.11%  "

Set-TextValue 24 4 "64.90"
$ws.Cells.Item(24, 5).Value = "  +2.45%  "

Set-TextValue 25 4 "1.75"
$ws.Cells.Item(25, 5).Value = "  +5.03%  "

Set-TextValue 26 4 "9.11"
$ws.Cells.Item(26, 5).Value = "  +8.66%  "

Set-TextValue 27 4 "577.23"
$ws.Cells.Item(27, 5).Value = "  +16.25%  "

$ws.Cells.Item(28, 4).Value = "2.542.09"
$ws.Cells.Item(28, 5).Value = "  +3.14%  "

$ws.Cells.Item(29, 5).Value = "  -0.12%  "

Set-TextValue 30 4 "8.38"
$ws.Cells.Item(30, 5).Value = "  +4.40%  "

$ws.Cells.Item(31, 4).Value = "0.0₃0932"
$ws.Cells.Item(31, 5).Value = "  +8.34%  "

$ws.Cells.Item(32, 5).Value = "  +5.77%  "

$ws.Cells.Item(33, 5).Value = "  +1.64%  "

$ws.Cells.Item(34, 5).Value = "  +4.06%  "

$ws.Cells.Item(35, 5).Value = "  +3.42%  "

$ws.Cells.Item(36, 5).Value = "  +0.01%  "

Set-TextValue 37 4 "5.66"
$ws.Cells.Item(37, 5).Value = "  +8.29%  "

Set-TextValue 38 4 "4.82"
$ws.Cells.Item(38, 5).Value = "  +4.16%  "

$ws.Cells.Item(39, 5).Value = "  +2.16%  "

$ws.Cells.Item(40, 2).Value = "EthereumClassic"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue 40 4 "18.72"
$ws.Cells.Item(40, 5).Value = "  +1.21%  "

$ws.Cells.Item(41, 2).Value = "Stacks"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue 41 4 "1.87"
$ws.Cells.Item(41, 5).Value = "  +2.91%  "

Set-TextValue 42 4 "148.58"
$ws.Cells.Item(42, 5).Value = "  +3.71%  "

$ws.Cells.Item(43, 5).Value = "  +0.05%  "

Set-TextValue 44 4 "41.67"
$ws.Cells.Item(44, 5).Value = "  +2.59%  "

Set-TextValue 45 4 "2.30"
$ws.Cells.Item(45, 5).Value = "  +13.53%  "

Set-TextValue 46 4 "150.85"
$ws.Cells.Item(46, 5).Value = "  +5.82%  "

$ws.Cells.Item(47, 5).Value = "  +2.07%  "

Set-TextValue 48 4 "0.0544"
$ws.Cells.Item(48, 5).Value = "  +5.59%  "

Set-TextValue 49 4 "20.32"
$ws.Cells.Item(49, 5).Value = "  +6.33%  "

Set-TextValue 50 4 "0.588"
$ws.Cells.Item(50, 5).Value = "  +3.68%  "

Set-TextValue 51 4 "0.0917"
$ws.Cells.Item(51, 5).Value = "  +1.94%  "
